$d = $word.ActiveDocument

# --- 1. First paragraph: update the bookmark placeholder text ---
$d.Content.Find.Execute("**ID__AFFARS_pgi_5301_topic_32__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_AFMC_PGI_5301_9001__ID**", 2) | Out-Null

$p = $d.Paragraphs(1)

# --- 2. Remove the now-orphaned trailing-space run that followed the bookmark text ---
$spaceRange = $d.Range($p.Range.End - 2, $p.Range.End - 1)
$spaceRange.Delete()

# --- 3. Paragraph formatting: indent + paragraph border spacing ---
$pf = $p.Range.ParagraphFormat
$pf.LeftIndent = 11.25
$pf.Borders.DistanceFromTop = 5
$pf.Borders.DistanceFromLeft = 5
$pf.Borders.DistanceFromBottom = 5
$pf.Borders.DistanceFromRight = 5

# --- 4. Bump the height of the "AFRL" row in the enterprise-contracting clearance-review table ---
$row = $d.Tables(3).Rows(3)
$row.Height = 84
